$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 36.307693
$ws.Range("I2").Value = 36.307693
$ws.Range("K2").Value = 36.307693
$ws.Range("M2").Value = 76.692307
$ws.Range("H4").Value = 5000
$ws.Range("I4").Value = 5000
$ws.Range("K4").Value = 5000
$ws.Range("M4").Value = -4886
$ws.Range("H5").Value = 77.69231000000001
$ws.Range("I5").Value = 72.84
$ws.Range("J5").Value = 199
$ws.Range("K5").Value = 72.84
$ws.Range("L5").Value = 199
$ws.Range("M5").Value = 42.16
$ws.Range("N5").Value = -429
$ws.Range("H17").Value = 470.875
$ws.Range("J17").Value = 554.087
$ws.Range("L17").Value = 1662.261
$ws.Range("N17").Value = -1998.261
$ws.Range("H68").Value = 72000
$ws.Range("J68").Value = 72000
$ws.Range("L68").Value = 72000
$ws.Range("N68").Value = -73498
$ws.Range("H71").Value = 72000
$ws.Range("J71").Value = 72000
$ws.Range("L71").Value = 216000
$ws.Range("N71").Value = -223488
$ws.Range("H82").Value = 200
$ws.Range("I82").Value = 200
$ws.Range("K82").Value = 600
$ws.Range("M82").Value = -194
$ws.Range("H85").Value = 200
$ws.Range("I85").Value = 200
$ws.Range("K85").Value = 600
$ws.Range("M85").Value = 804
$ws.Range("H98").Value = 38465484
$ws.Range("I98").Value = 40003668
$ws.Range("K98").Value = 40003668
$ws.Range("M98").Value = -40002170
$ws.Range("H122").Value = 38465484
$ws.Range("I122").Value = 40003668
$ws.Range("K122").Value = 120011004
$ws.Range("M122").Value = -120008554

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 64.75
$ws.Range("I5").Value = 64.75
$ws.Range("K5").Value = 64.75
$ws.Range("M5").Value = 47.25
$ws.Range("H61").Value = 2149.7646
$ws.Range("I61").Value = 2168.5312
$ws.Range("K61").Value = 2168.5312
$ws.Range("M61").Value = -1956.5312
$ws.Range("H74").Value = 50181.734
$ws.Range("I74").Value = 69352.914
$ws.Range("K74").Value = 69352.914
$ws.Range("M74").Value = -68478.914
$ws.Range("H77").Value = 50181.734
$ws.Range("I77").Value = 69352.914
$ws.Range("K77").Value = 346764.57
$ws.Range("M77").Value = -342396.57
$ws.Range("H132").Value = 10664.521
$ws.Range("I132").Value = 11003.4
$ws.Range("K132").Value = 33010.2
$ws.Range("M132").Value = -30480.2
$ws.Range("H136").Value = 2149.7646
$ws.Range("I136").Value = 2168.5312
$ws.Range("K136").Value = 6505.5936
$ws.Range("M136").Value = -3955.5936

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 64.75
$ws.Range("I4").Value = 64.75
$ws.Range("K4").Value = 64.75
$ws.Range("M4").Value = 50.25
$ws.Range("H99").Value = 4135851.5
$ws.Range("I99").Value = 2393.4
$ws.Range("J99").Value = 7580400
$ws.Range("K99").Value = 2393.4
$ws.Range("L99").Value = 7580400
$ws.Range("M99").Value = -895.4000000000001
$ws.Range("N99").Value = -7583396

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 7500
$ws.Range("I3").Value = 7500
$ws.Range("K3").Value = 7500
$ws.Range("M3").Value = -7387
$ws.Range("H58").Value = 8934321
$ws.Range("I58").Value = 21741882
$ws.Range("J58").Value = 7838.9395
$ws.Range("K58").Value = 21741882
$ws.Range("L58").Value = 7838.9395
$ws.Range("M58").Value = -21741679
$ws.Range("N58").Value = -8244.9395
$ws.Range("H134").Value = 8793.565000000001
$ws.Range("I134").Value = 9376.700000000001
$ws.Range("J134").Value = 8345
$ws.Range("K134").Value = 28130.1
$ws.Range("L134").Value = 25035
$ws.Range("M134").Value = -25595.1
$ws.Range("N134").Value = -30105
$ws.Range("H136").Value = 8934321
$ws.Range("I136").Value = 21741882
$ws.Range("J136").Value = 7838.9395
$ws.Range("K136").Value = 65225646
$ws.Range("L136").Value = 23516.8185
$ws.Range("M136").Value = -65223096
$ws.Range("N136").Value = -28616.8185

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8335204
$ws.Range("J4").Value = 32577.5
$ws.Range("L4").Value = 97732.5
$ws.Range("N4").Value = -97956.5
$ws.Range("H7").Value = 1740
$ws.Range("I7").Value = 266.66666
$ws.Range("K7").Value = 799.9999799999999
$ws.Range("M7").Value = -687.9999799999999
$ws.Range("H23").Value = 292.25
$ws.Range("J23").Value = 293.6
$ws.Range("L23").Value = 880.8000000000001
$ws.Range("N23").Value = -1350.8
$ws.Range("H33").Value = 25641242
$ws.Range("I33").Value = 55555630
$ws.Range("K33").Value = 333333780
$ws.Range("M33").Value = -333333497
$ws.Range("H44").Value = 1720
$ws.Range("I44").Value = 300
$ws.Range("J44").Value = 2666.6667
$ws.Range("K44").Value = 900
$ws.Range("L44").Value = 8000.000100000001
$ws.Range("M44").Value = -502
$ws.Range("N44").Value = -8796.000100000001
$ws.Range("H134").Value = 82041.46000000001
$ws.Range("I134").Value = 127067.375
$ws.Range("K134").Value = 381202.125
$ws.Range("M134").Value = -376132.125
$ws.Range("H138").Value = 66269.69
$ws.Range("I138").Value = 86817.914
$ws.Range("K138").Value = 260453.742
$ws.Range("M138").Value = -255313.742
$ws.Range("H139").Value = 38603.332
$ws.Range("I139").Value = 68773.47
$ws.Range("J139").Value = 8433.200000000001
$ws.Range("K139").Value = 206320.41
$ws.Range("L139").Value = 25299.6
$ws.Range("M139").Value = -201180.41
$ws.Range("N139").Value = -35579.60000000001
$ws.Range("H140").Value = 167740.67
$ws.Range("I140").Value = 167740.67
$ws.Range("K140").Value = 503222.01
$ws.Range("M140").Value = -498042.01

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 8368.286
$ws.Range("I132").Value = 4505.75
$ws.Range("K132").Value = 13517.25
$ws.Range("M132").Value = -10987.25
$ws.Range("H136").Value = 60055
$ws.Range("J136").Value = 60055
$ws.Range("L136").Value = 180165
$ws.Range("N136").Value = -185265

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7475.5
$ws.Range("I61").Value = 5752
$ws.Range("K61").Value = 5752
$ws.Range("M61").Value = -5550
$ws.Range("H69").Value = 44998.5
$ws.Range("J69").Value = 44998.5
$ws.Range("L69").Value = 44998.5
$ws.Range("N69").Value = -46620.5
$ws.Range("H72").Value = 44998.5
$ws.Range("J72").Value = 44998.5
$ws.Range("L72").Value = 134995.5
$ws.Range("N72").Value = -143107.5
$ws.Range("H82").Value = 1085467.8
$ws.Range("I82").Value = 1566987.5
$ws.Range("K82").Value = 1566987.5
$ws.Range("M82").Value = -1566626.5
$ws.Range("H85").Value = 1085467.8
$ws.Range("I85").Value = 1566987.5
$ws.Range("K85").Value = 1566987.5
$ws.Range("M85").Value = -1565739.5
$ws.Range("H100").Value = 3594.6316
$ws.Range("I100").Value = 2962.8333
$ws.Range("J100").Value = 3886.2307
$ws.Range("K100").Value = 2962.8333
$ws.Range("L100").Value = 3886.2307
$ws.Range("M100").Value = -2421.8333
$ws.Range("N100").Value = -4968.2307
$ws.Range("H113").Value = 7475.5
$ws.Range("I113").Value = 5752
$ws.Range("K113").Value = 5752
$ws.Range("M113").Value = -3582

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 174999.72
$ws.Range("I15").Value = 174999.72
$ws.Range("K15").Value = 174999.72
$ws.Range("M15").Value = -174711.72
$ws.Range("H43").Value = 19200
$ws.Range("I43").Value = 12000
$ws.Range("K43").Value = 12000
$ws.Range("M43").Value = -11851
$ws.Range("H96").Value = 1923.2858
$ws.Range("I96").Value = 1997.25
$ws.Range("K96").Value = 1997.25
$ws.Range("M96").Value = -624.25
$ws.Range("H122").Value = 140989.38
$ws.Range("I122").Value = 236628.94
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 709886.8200000001
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -707436.8200000001
$ws.Range("N122").Value = -21400
$ws.Range("H132").Value = 14291572
$ws.Range("I132").Value = 22733430
$ws.Range("K132").Value = 68200290
$ws.Range("M132").Value = -68197760
